# This script applies proofing-mark (w:proofErr) annotations that Word
# inserts after a spell/grammar check pass, splitting a few runs so the
# flagged words/phrases sit in their own <w:r>, and also removes the
# (now redundant) " au niveau de l'" trailing text in one cell.
#
# Because the host object model does not auto-generate <w:proofErr/>
# markers, each affected paragraph is rebuilt explicitly via
# Range.InsertXML with a single-part WordprocessingML package, which
# lets us place <w:proofErr> elements between/around <w:r> runs exactly
# as Word's proofing pass would.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
  '<w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Replace-Paragraph {
    param(
        [string]$FindText,
        [string]$ParagraphXml
    )
    $r = $d.Content
    $found = $r.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $FindText"
        return
    }
    $p = $r.Paragraphs(1).Range
    $xml = $pkgHeader + $ParagraphXml + $pkgFooter
    $ret = $p.InsertXML($xml)
}

# --- "aaaa-mm-dd" date-format hints: 5 occurrences, each becomes
#     spellStart/gramStart "aaaa" + spellEnd/gramEnd, then "-mm-dd" ---

$find = "aaaa-mm-dd"

$xml1 =
  '<w:p w14:paraId="2067C7AA" w14:textId="6168BA67" w:rsidR="00B744A9" w:rsidRDefault="0086049A" w:rsidP="003E75EE">' +
  '<w:pPr><w:ind w:left="0" w:firstLine="0"/></w:pPr>' +
  '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>aaaa</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t>-mm-dd</w:t></w:r>' +
  '</w:p>'
Replace-Paragraph $find $xml1

$xml2 =
  '<w:p w14:paraId="11D499CE" w14:textId="3BC1924A" w:rsidR="0086049A" w:rsidRDefault="00A50B5B" w:rsidP="003E75EE">' +
  '<w:pPr><w:ind w:left="0" w:firstLine="0"/></w:pPr>' +
  '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>aaaa</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t>-mm-dd</w:t></w:r>' +
  '</w:p>'
Replace-Paragraph $find $xml2

$xml3 =
  '<w:p w14:paraId="00D92F97" w14:textId="50B9BB06" w:rsidR="0086049A" w:rsidRDefault="00A50B5B" w:rsidP="003E75EE">' +
  '<w:pPr><w:ind w:left="0" w:firstLine="0"/></w:pPr>' +
  '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>aaaa</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t>-mm-dd</w:t></w:r>' +
  '</w:p>'
Replace-Paragraph $find $xml3

$xml4 =
  '<w:p w14:paraId="47F7FC1C" w14:textId="0566D376" w:rsidR="00B141E8" w:rsidRDefault="00B141E8" w:rsidP="003E75EE">' +
  '<w:pPr><w:ind w:left="0" w:firstLine="0"/></w:pPr>' +
  '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>aaaa</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t>-mm-dd</w:t></w:r>' +
  '</w:p>'
Replace-Paragraph $find $xml4

$xml5 =
  '<w:p w14:paraId="14A9A011" w14:textId="57ADA90D" w:rsidR="00B141E8" w:rsidRDefault="00B141E8" w:rsidP="003B5E3C">' +
  '<w:pPr><w:ind w:left="0" w:firstLine="0"/></w:pPr>' +
  '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>aaaa</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t>-mm-dd</w:t></w:r>' +
  '</w:p>'
Replace-Paragraph $find $xml5

# --- "« male » ou « female »" : wrap "male" (gramStart/End) and
#     "female" (spellStart/End) ---

$find = "female"
$xml =
  '<w:p w14:paraId="7F9CB58D" w14:textId="287816F0" w:rsidR="006F4789" w:rsidRDefault="00C14E20" w:rsidP="003E75EE">' +
  '<w:pPr><w:ind w:left="0" w:firstLine="0"/></w:pPr>' +
  '<w:r><w:t>« </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>m</w:t></w:r>' +
  '<w:r w:rsidR="0086049A"><w:t>ale</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t> »</w:t></w:r>' +
  '<w:r w:rsidR="0086049A"><w:t xml:space="preserve"> ou </w:t></w:r>' +
  '<w:r><w:t>« </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r w:rsidR="0086049A"><w:t>female</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t> »</w:t></w:r>' +
  '</w:p>'
Replace-Paragraph $find $xml

# --- "Adresse email valide et unique au niveau de l'" -> drop the
#     trailing "au niveau de l'" text, keeping a lone space run ---

$apos = [char]0x2019
$find = "au niveau de l" + $apos
$xml =
  '<w:p w14:paraId="3AA272EC" w14:textId="062EEFC7" w:rsidR="00D0760A" w:rsidRDefault="00D0760A" w:rsidP="003E75EE">' +
  '<w:pPr><w:ind w:left="0" w:firstLine="0"/></w:pPr>' +
  '<w:r><w:t>Adresse email valide</w:t></w:r>' +
  '<w:r w:rsidR="000C282A"><w:t xml:space="preserve"> et unique</w:t></w:r>' +
  '<w:r w:rsidR="00136BA7"><w:t xml:space="preserve"> </w:t></w:r>' +
  '</w:p>'
Replace-Paragraph $find $xml

# --- "Commence par « it » et se termine par un chiffre" : wrap "it"
#     with spellStart/spellEnd ---

$find = "Commence par"
$xml =
  '<w:p w14:paraId="5289996C" w14:textId="4CC46E7E" w:rsidR="00B141E8" w:rsidRDefault="006C50AC" w:rsidP="00B64991">' +
  '<w:pPr><w:ind w:left="0" w:firstLine="0"/></w:pPr>' +
  '<w:r><w:t>Commence par « </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>it</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t> » et se termine par un chiffre</w:t></w:r>' +
  '</w:p>'
Replace-Paragraph $find $xml

# --- "open, resolved, closed" : wrap each status word ---

$find = "open, resolved, closed"
$xml =
  '<w:p w14:paraId="3A23E287" w14:textId="26FE1253" w:rsidR="004D3D19" w:rsidRDefault="00B3673D" w:rsidP="00B64991">' +
  '<w:pPr><w:ind w:left="0" w:firstLine="0"/></w:pPr>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>open</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>resolved</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>closed</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>'
Replace-Paragraph $find $xml

# --- "S'inspirer de la vue de connexion de openclassroom" : wrap
#     "openclassroom" with spellStart/spellEnd (no w:pPr on this one) ---

$find = "openclassroom"
$xml =
  '<w:p w14:paraId="50B7C782" w14:textId="2DB2AF92" w:rsidR="002849AD" w:rsidRDefault="002849AD" w:rsidP="002849AD">' +
  '<w:r><w:t xml:space="preserve">S’inspirer de la vue de connexion de </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>openclassroom</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>'
Replace-Paragraph $find $xml
